# Update localization status from "Ready for handoff" to "In Translation"
# for the two files that have moved into translation:
#   0c2784e3-dc6c-4956-8952-0bc931416b6c.md
#   0d6fefdb-ca00-4562-b000-12caef14fceb.md
# (831e25f3-a181-4e16-9403-096be0873547.md stays "Ready for handoff")

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: zh-cn (E) and de-de (F) status columns, rows 3 & 4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- "zh-cn" detail sheet: Status column (C), rows 3 & 4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- "de-de" detail sheet: Status column (C), rows 3 & 4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
